$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9511678814888
$ws.Range("B1").Value = 1.741506457328796
$ws.Range("C1").Value = 6.592920780181885
$ws.Range("D1").Value = 3.345165967941284
$ws.Range("E1").Value = 1.507930874824524
